# Georgia -> DejaVu Sans for repeatable layout test
#
# Every style whose own rPr explicitly pins "Georgia" as the ascii/hAnsi
# font gets switched to "DejaVu Sans" (eastAsiaTheme/cstheme and any other
# run-property stay untouched; Word's Font.Name setter only rewrites the
# ascii/hAnsi font slots).

$d = $word.ActiveDocument

$styleNames = @(
    "Normal",
    "Heading 1 Char",
    "Heading 2 Char",
    "Header Char",
    "Footer Char",
    "Title Char",
    "No Spacing",
    "Subtitle Char",
    "Body Text Char",
    "No Spacing Char",
    "Comment Char"
)

foreach ($name in $styleNames) {
    $style = $d.Styles.Item($name)
    $style.Font.Name = "DejaVu Sans"
}

Write-Output "Updated $($styleNames.Count) styles from Georgia to DejaVu Sans"
